$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing text storage (avoids Excel auto-numeric
# conversion for plain decimal-looking strings) and keeping the default "Normal"
# cell style (no visible NumberFormat change).
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# --- Rows 2-46: update Price (D) and Volume(1h) (E) columns ---
$ws.Range("D2").Value = "26.464.41"
$ws.Range("E2").Value = "  -1.09%  "

$ws.Range("D3").Value = "1.625.25"
$ws.Range("E3").Value = "  -0.84%  "

$ws.Range("E4").Value = "  +0.27%  "

Set-TextValue "D5" "212.60"
$ws.Range("E5").Value = "  -0.20%  "

Set-TextValue "D6" "0.499"
$ws.Range("E6").Value = "  +0.91%  "

$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("E8").Value = "  -0.54%  "

$ws.Range("E9").Value = "  -2.15%  "

Set-TextValue "D10" "18.72"
$ws.Range("E10").Value = "  -1.96%  "

$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("D12").Value = "1.851.09"
$ws.Range("E12").Value = "  -0.83%  "

$ws.Range("D13").Value = "1.633.48"
$ws.Range("E13").Value = "  -0.76%  "

$ws.Range("E14").Value = "  +1.03%  "

$ws.Range("E15").Value = "  -1.02%  "

Set-TextValue "D16" "64.81"
$ws.Range("E16").Value = "  +2.56%  "

$ws.Range("D17").Value = "26.475.31"
$ws.Range("E17").Value = "  -1.00%  "

$ws.Range("E18").Value = "  -0.73%  "

Set-TextValue "D19" "213.34"
$ws.Range("E19").Value = "  +1.98%  "

$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("E21").Value = "  -0.90%  "

$ws.Range("E22").Value = "  +1.19%  "

Set-TextValue "D23" "9.26"
$ws.Range("E23").Value = "  -1.65%  "

Set-TextValue "D24" "2.03"
$ws.Range("E24").Value = "  +5.18%  "

Set-TextValue "D25" "148.52"
$ws.Range("E25").Value = "  +1.50%  "

Set-TextValue "D26" "1.00"
$ws.Range("E26").Value = "  +0.33%  "

$ws.Range("E27").Value = "  -1.10%  "

$ws.Range("E28").Value = "  +0.76%  "

Set-TextValue "D29" "15.47"
$ws.Range("E29").Value = "  +0.43%  "

Set-TextValue "D30" "0.0507"
$ws.Range("E30").Value = "  -1.84%  "

$ws.Range("E31").Value = "  -1.05%  "

$ws.Range("E32").Value = "  +2.31%  "

$ws.Range("E33").Value = "  -1.24%  "

$ws.Range("E34").Value = "  -0.59%  "

$ws.Range("D35").Value = "1.221.93"
$ws.Range("E35").Value = "  +4.44%  "

$ws.Range("E36").Value = "  -1.10%  "

Set-TextValue "D37" "0.0172"
$ws.Range("E37").Value = "  +3.06%  "

$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("E39").Value = "  -2.30%  "

Set-TextValue "D40" "0.504"
$ws.Range("E40").Value = "  -0.09%  "

Set-TextValue "D41" "2.28"
$ws.Range("E41").Value = "  -1.79%  "

$ws.Range("E42").Value = "  -0.50%  "

$ws.Range("E43").Value = "  -0.80%  "

$ws.Range("D44").Value = "1.760.57"
$ws.Range("E44").Value = "  -0.84%  "

Set-TextValue "D45" "92.70"
$ws.Range("E45").Value = "  +0.26%  "

Set-TextValue "D46" "1.57"
$ws.Range("E46").Value = "  +0.54%  "

# --- Rows 47-51: coin list entries shifted up by one; update Coin, Link, Price, Volume(1h) ---
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D47" "54.64"
$ws.Range("E47").Value = "  -0.17%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D48" "0.0510"
$ws.Range("E48").Value = "  -0.53%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D49" "0.406"
$ws.Range("E49").Value = "  -0.97%  "

$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextValue "D50" "1.00"
$ws.Range("E50").Value = "  +0.29%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "7.46"
$ws.Range("E51").Value = "  -0.78%  "
